{"js": "// Locate the (only) table in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Row 1 (0-indexed), Column 1 (0-indexed) -> \"NAME TASK\" cell for \"Sprint 5\".\n// Replace the numbered-list paragraph (\"Pembuatan Desain dari fitur Login\")\n// with a plain paragraph made of three runs: \"1. \" / \"Pembuatan Desain dari\n// fitur \" / \"Bookmark\".\nconst taskCell = table.getCell(1, 1);\nconst taskOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:t xml:space=\"preserve\">1. </w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">Pembuatan Desain dari fitur </w:t></w:r>\n<w:r><w:t>Bookmark</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\ntaskCell.body.insertOoxml(taskOoxml, \"Replace\");\n\n// Row 1, Column 2 -> the \"1 Day\" estimate cell. Drop the \"_GoBack\" bookmark\n// that used to sit between the \"1\" and \" Day\" runs (it moves to the new\n// row below), keeping just the two text runs.\nconst estimateCell = table.getCell(1, 2);\nconst estimateOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:t>1</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\"> Day</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\nestimateCell.body.insertOoxml(estimateOoxml, \"Replace\");\nawait context.sync();\n\n// Append a brand-new \"Sprint 6\" row with placeholder text, then fill in the\n// task-description cell (three runs) and the estimate cell (single run +\n// a freshly created \"_GoBack\" bookmark) with precise OOXML.\ntable.addRows(\"End\", 1, [[\"Sprint 6\", \"\", \"\"]]);\nawait context.sync();\n\nconst newTaskCell = table.getCell(2, 1);\nconst newTaskOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:t xml:space=\"preserve\">2. </w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">Pembuatan Desain dari fitur </w:t></w:r>\n<w:r><w:t>Profil</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\nnewTaskCell.body.insertOoxml(newTaskOoxml, \"Replace\");\n\nconst newEstimateCell = table.getCell(2, 2);\nconst newEstimateOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p>\n<w:r><w:t>1 Day</w:t></w:r>\n<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n<w:bookmarkEnd w:id=\"0\"/>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\nnewEstimateCell.body.insertOoxml(newEstimateOoxml, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Word COM interop script.\n# $word.ActiveDocument is the open document ($d is a convenience alias).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# --- Row \"Sprint 5\" (table row 2) ---------------------------------------\n\n# NAME TASK cell: drop the numbered-list paragraph formatting and collapse\n# the five spell-checked fragments into three plain runs:\n# \"1. \" / \"Pembuatan Desain dari fitur \" / \"Bookmark\".\n$taskXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">1. </w:t></w:r><w:r><w:t xml:space=\"preserve\">Pembuatan Desain dari fitur </w:t></w:r><w:r><w:t>Bookmark</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$t.Cell(2, 2).Range.InsertXML($taskXml)\n\n# ESTIMATE PROCESSING TIME cell: remove the \"_GoBack\" bookmark that used to\n# sit between \"1\" and \" Day\" (it gets re-created on the new row below).\n$estimateXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>1</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Day</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$t.Cell(2, 3).Range.InsertXML($estimateXml)\n\n# --- New \"Sprint 6\" row (table row 3) -----------------------------------\n\n$t.Rows.Add() | Out-Null\n\n$sprintXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>Sprint 6</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$t.Cell(3, 1).Range.InsertXML($sprintXml)\n\n$newTaskXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">2. </w:t></w:r><w:r><w:t xml:space=\"preserve\">Pembuatan Desain dari fitur </w:t></w:r><w:r><w:t>Profil</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$t.Cell(3, 2).Range.InsertXML($newTaskXml)\n\n$newEstimateXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>1 Day</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$t.Cell(3, 3).Range.InsertXML($newEstimateXml)\n"}
